$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1519
$ws.Cells.Item(32, 10).Value = 1523.75
$ws.Cells.Item(32, 12).Value = 1523.75
$ws.Cells.Item(32, 14).Value = -2175.75
$ws.Cells.Item(74, 8).Value = 3000
$ws.Cells.Item(74, 9).Value = 3000
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 3000
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -2064
$ws.Cells.Item(74, 14).Value = ""
$ws.Cells.Item(77, 8).Value = 3000
$ws.Cells.Item(77, 9).Value = 3000
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 15000
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -10320
$ws.Cells.Item(77, 14).Value = ""
$ws.Cells.Item(100, 8).Value = 2381
$ws.Cells.Item(100, 9).Value = 2330.4285
$ws.Cells.Item(100, 10).Value = 2499
$ws.Cells.Item(100, 11).Value = 2330.4285
$ws.Cells.Item(100, 12).Value = 2499
$ws.Cells.Item(100, 13).Value = -1789.4285
$ws.Cells.Item(100, 14).Value = -3581
$ws.Cells.Item(138, 8).Value = 2198.43
$ws.Cells.Item(138, 9).Value = 1294.6666
$ws.Cells.Item(138, 10).Value = 2256.117
$ws.Cells.Item(138, 11).Value = 3883.9998
$ws.Cells.Item(138, 12).Value = 6768.351000000001
$ws.Cells.Item(138, 13).Value = 1256.0002
$ws.Cells.Item(138, 14).Value = -17048.351
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2355.4814
$ws.Cells.Item(32, 9).Value = 2572.4583
$ws.Cells.Item(32, 11).Value = 2572.4583
$ws.Cells.Item(32, 13).Value = -2285.4583
$ws.Cells.Item(45, 8).Value = 2140.0908
$ws.Cells.Item(45, 9).Value = 2317.625
$ws.Cells.Item(45, 10).Value = 1666.6666
$ws.Cells.Item(45, 11).Value = 2317.625
$ws.Cells.Item(45, 12).Value = 1666.6666
$ws.Cells.Item(45, 13).Value = -1940.625
$ws.Cells.Item(45, 14).Value = -2420.6666
$ws.Cells.Item(61, 8).Value = 1323.7222
$ws.Cells.Item(61, 9).Value = 1225.1177
$ws.Cells.Item(61, 10).Value = 3000
$ws.Cells.Item(61, 11).Value = 1225.1177
$ws.Cells.Item(61, 12).Value = 3000
$ws.Cells.Item(61, 13).Value = -1013.1177
$ws.Cells.Item(61, 14).Value = -3424
$ws.Cells.Item(63, 8).Value = 2250
$ws.Cells.Item(63, 9).Value = 1950
$ws.Cells.Item(63, 10).Value = 2550
$ws.Cells.Item(63, 11).Value = 1950
$ws.Cells.Item(63, 12).Value = 2550
$ws.Cells.Item(63, 13).Value = -1264
$ws.Cells.Item(63, 14).Value = -3922
$ws.Cells.Item(66, 8).Value = 2250
$ws.Cells.Item(66, 9).Value = 1950
$ws.Cells.Item(66, 10).Value = 2550
$ws.Cells.Item(66, 11).Value = 9750
$ws.Cells.Item(66, 12).Value = 12750
$ws.Cells.Item(66, 13).Value = -6318
$ws.Cells.Item(66, 14).Value = -19614
$ws.Cells.Item(102, 8).Value = 27794896
$ws.Cells.Item(102, 9).Value = 27794896
$ws.Cells.Item(102, 11).Value = 27794896
$ws.Cells.Item(102, 13).Value = -27793274
$ws.Cells.Item(122, 8).Value = 1764.8
$ws.Cells.Item(122, 9).Value = 1778
$ws.Cells.Item(122, 10).Value = 1712
$ws.Cells.Item(122, 11).Value = 5334
$ws.Cells.Item(122, 12).Value = 5136
$ws.Cells.Item(122, 13).Value = -2884
$ws.Cells.Item(122, 14).Value = -10036
$ws.Cells.Item(136, 8).Value = 1323.7222
$ws.Cells.Item(136, 9).Value = 1225.1177
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 3675.3531
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = -1125.3531
$ws.Cells.Item(136, 14).Value = -14100
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 90911480
$ws.Cells.Item(105, 9).Value = 111113250
$ws.Cells.Item(105, 11).Value = 111113250
$ws.Cells.Item(105, 13).Value = -111111503
$ws.Cells.Item(107, 8).Value = 1538.6428
$ws.Cells.Item(107, 9).Value = 1192.8
$ws.Cells.Item(107, 10).Value = 2403.25
$ws.Cells.Item(107, 11).Value = 1192.8
$ws.Cells.Item(107, 12).Value = 2403.25
$ws.Cells.Item(107, 13).Value = 727.2
$ws.Cells.Item(107, 14).Value = -6243.25
$ws.Cells.Item(134, 8).Value = 9175.629999999999
$ws.Cells.Item(134, 9).Value = 6378.2383
$ws.Cells.Item(134, 11).Value = 19134.7149
$ws.Cells.Item(134, 13).Value = -16599.7149
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 564
$ws.Cells.Item(105, 9).Value = 410
$ws.Cells.Item(105, 10).Value = 666.6667
$ws.Cells.Item(105, 11).Value = 410
$ws.Cells.Item(105, 12).Value = 666.6667
$ws.Cells.Item(105, 13).Value = 1337
$ws.Cells.Item(105, 14).Value = -4160.6667
$ws.Cells.Item(132, 8).Value = 8495.058999999999
$ws.Cells.Item(132, 9).Value = 15657.714
$ws.Cells.Item(132, 10).Value = 3481.2
$ws.Cells.Item(132, 11).Value = 46973.142
$ws.Cells.Item(132, 12).Value = 10443.6
$ws.Cells.Item(132, 13).Value = -44443.142
$ws.Cells.Item(132, 14).Value = -15503.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 618.5
$ws.Cells.Item(12, 10).Value = 634.8182
$ws.Cells.Item(12, 12).Value = 1904.4546
$ws.Cells.Item(12, 14).Value = -2250.4546
$ws.Cells.Item(107, 8).Value = 4908.5415
$ws.Cells.Item(107, 9).Value = 603.61536
$ws.Cells.Item(107, 11).Value = 1810.84608
$ws.Cells.Item(107, 13).Value = 109.15392
$ws.Cells.Item(122, 8).Value = 709.0526
$ws.Cells.Item(122, 10).Value = 756.7059
$ws.Cells.Item(122, 12).Value = 6810.3531
$ws.Cells.Item(122, 14).Value = -11710.3531
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 13185.714
$ws.Cells.Item(5, 10).Value = 13185.714
$ws.Cells.Item(5, 12).Value = 13185.714
$ws.Cells.Item(5, 14).Value = -13409.714
$ws.Cells.Item(126, 8).Value = 2294.1538
$ws.Cells.Item(126, 9).Value = 1769.3334
$ws.Cells.Item(126, 11).Value = 5308.0002
$ws.Cells.Item(126, 13).Value = -2838.0002
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 436111.12
$ws.Cells.Item(2, 10).Value = 372222.22
$ws.Cells.Item(2, 12).Value = 372222.22
$ws.Cells.Item(2, 14).Value = -372446.22
$ws.Cells.Item(16, 8).Value = 1055.8636
$ws.Cells.Item(16, 9).Value = 979.2778
$ws.Cells.Item(16, 10).Value = 1400.5
$ws.Cells.Item(16, 11).Value = 979.2778
$ws.Cells.Item(16, 12).Value = 1400.5
$ws.Cells.Item(16, 13).Value = -809.2778
$ws.Cells.Item(16, 14).Value = -1740.5
$ws.Cells.Item(40, 8).Value = 3003
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 3003
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 3003
$ws.Cells.Item(40, 13).Value = ""
$ws.Cells.Item(40, 14).Value = -3275
$ws.Cells.Item(61, 8).Value = 3141.4
$ws.Cells.Item(61, 9).Value = 2400.6667
$ws.Cells.Item(61, 10).Value = 4252.5
$ws.Cells.Item(61, 11).Value = 2400.6667
$ws.Cells.Item(61, 12).Value = 4252.5
$ws.Cells.Item(61, 13).Value = -2198.6667
$ws.Cells.Item(61, 14).Value = -4656.5
$ws.Cells.Item(100, 8).Value = 1620.3
$ws.Cells.Item(100, 9).Value = 1412.875
$ws.Cells.Item(100, 10).Value = 2450
$ws.Cells.Item(100, 11).Value = 1412.875
$ws.Cells.Item(100, 12).Value = 2450
$ws.Cells.Item(100, 13).Value = -871.875
$ws.Cells.Item(100, 14).Value = -3532
$ws.Cells.Item(113, 8).Value = 3141.4
$ws.Cells.Item(113, 9).Value = 2400.6667
$ws.Cells.Item(113, 10).Value = 4252.5
$ws.Cells.Item(113, 11).Value = 2400.6667
$ws.Cells.Item(113, 12).Value = 4252.5
$ws.Cells.Item(113, 13).Value = -230.6667000000002
$ws.Cells.Item(113, 14).Value = -8592.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 636.6429000000001
$ws.Cells.Item(113, 9).Value = 326.3
$ws.Cells.Item(113, 10).Value = 1412.5
$ws.Cells.Item(113, 11).Value = 978.9000000000001
$ws.Cells.Item(113, 12).Value = 4237.5
$ws.Cells.Item(113, 13).Value = 1191.1
$ws.Cells.Item(113, 14).Value = -8577.5
$ws.Cells.Item(136, 8).Value = 1551
$ws.Cells.Item(136, 9).Value = 1337.2727
$ws.Cells.Item(136, 11).Value = 4011.8181
$ws.Cells.Item(136, 13).Value = -1461.8181

Write-Host "Applied all changes"